$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 19
$ws.Range("A3").Value = 20
$ws.Range("A4").Value = 21
$ws.Range("A5").Value = 22
$ws.Range("A6").Value = 23

$ws.Range("H8").Select()
